# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF -------------------------
# Copy the formatting of the last existing header cell (AC1 - bold, border,
# centered) onto the three new header cells, then overwrite their text.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-61): season record for every player --------------------
$wins = 75
$losses = 87
$ties = 0

for ($row = 2; $row -le 61; $row++) {
    $ws.Cells.Item($row, 30).Value = $wins    # column AD
    $ws.Cells.Item($row, 31).Value = $losses  # column AE
    $ws.Cells.Item($row, 32).Value = $ties    # column AF
}
